$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings, several of which look like plain decimal
# numbers (e.g. "209.97", "0.870"). Assigning those directly to .Value
# would make Excel coerce them to real numbers (and e.g. drop the
# trailing zero in "0.870"). Temporarily force the column to Text format
# while writing the values, then restore the default "Normal" style so
# the cells end up exactly as before (general/no explicit style) but
# still holding literal text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.712.02'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '1.566.41'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  -0.43%  '

$ws.Range("D5").Value = '209.97'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").Value = '0.489'
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.48%  '

$ws.Range("D8").Value = '25.11'
$ws.Range("E8").Value = '  +5.98%  '

$ws.Range("D9").Value = '0.245'
$ws.Range("E9").Value = '  +0.42%  '

$ws.Range("D10").Value = '0.0586'
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").Value = '0.0895'
$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").Value = '1.790.00'
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").Value = '1.573.03'
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").Value = '28.692.90'
$ws.Range("E14").Value = '  +1.34%  '

$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").Value = '61.44'
$ws.Range("E17").Value = '  +0.66%  '

$ws.Range("D18").Value = '230.08'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").Value = '7.36'
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("D20").Value = '0.0₃0680'
$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").Value = '9.03'
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  +1.55%  '

$ws.Range("D25").Value = '151.51'
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("D26").Value = '14.81'
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("E27").Value = '  +0.40%  '

$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("D29").Value = '6.24'

$ws.Range("E30").Value = '  -3.66%  '

$ws.Range("E31").Value = '  -2.40%  '

$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("D33").Value = '1.395.15'
$ws.Range("E33").Value = '  +1.19%  '

$ws.Range("E34").Value = '  -2.81%  '

$ws.Range("E35").Value = '  -3.72%  '

$ws.Range("D36").Value = '1.48'
$ws.Range("E36").Value = '  -1.11%  '

$ws.Range("E37").Value = '  +2.24%  '

$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("D40").Value = '1.97'
$ws.Range("E40").Value = '  +2.15%  '

$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.38%  '

$ws.Range("D43").Value = '0.773'
$ws.Range("E43").Value = '  -1.14%  '

$ws.Range("E44").Value = '  -3.64%  '

$ws.Range("D45").Value = '63.97'
$ws.Range("E45").Value = '  +2.98%  '

$ws.Range("D46").Value = '5.24'
$ws.Range("E46").Value = '  -1.58%  '

$ws.Range("D47").Value = '1.701.83'
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").Value = '0.870'
$ws.Range("E48").Value = '  -5.07%  '

$ws.Range("D49").Value = '85.04'
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("D50").Value = '43.07'
$ws.Range("E50").Value = '  +6.04%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0512'
$ws.Range("E51").Value = '  -0.34%  '

# Restore the column's original (default/general) style now that every
# text value is locked in, so only the cell contents changed.
$ws.Range("D2:D51").Style = "Normal"
